$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple per-row price / volume(1h) updates ---
$ws.Range("D2").Value = "29.417.15"
$ws.Range("E2").Value = "  -0.02%  "

$ws.Range("D3").Value = "1.916.70"
$ws.Range("E3").Value = "  +0.78%  "

$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.52%  "

$ws.Range("D5").Value = "324.83"

$ws.Range("E6").Value = "  +0.37%  "

$ws.Range("D7").Value = "0.4817"
$ws.Range("E7").Value = "  +0.37%  "

$ws.Range("D8").Value = "0.4068"
$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").Value = "0.08245"
$ws.Range("E9").Value = "  +2.17%  "

$ws.Range("D10").Value = "1.012"
$ws.Range("E10").Value = "  +0.89%  "

$ws.Range("E11").Value = "  +0.17%  "

$ws.Range("D12").Value = "1.925.66"
$ws.Range("E12").Value = "  -0.63%  "

$ws.Range("D13").Value = "6.065"
$ws.Range("E13").Value = "  +1.95%  "

$ws.Range("E14").Value = "  +2.35%  "

$ws.Range("D15").Value = "91.59"
$ws.Range("E15").Value = "  +1.89%  "

$ws.Range("D16").Value = "0.06841"
$ws.Range("E16").Value = "  +2.02%  "

$ws.Range("E17").Value = "  +0.36%  "

$ws.Range("D18").Value = "0.00001039"
$ws.Range("E18").Value = "  +0.70%  "

$ws.Range("D19").Value = "17.62"
$ws.Range("E19").Value = "  -0.15%  "

$ws.Range("E20").Value = "  +0.42%  "

$ws.Range("D21").Value = "29.421.04"
$ws.Range("E21").Value = "  -0.06%  "

$ws.Range("E22").Value = "  +2.04%  "

$ws.Range("D23").Value = "11.75"
$ws.Range("E23").Value = "  -0.13%  "

$ws.Range("D24").Value = "2.177"
$ws.Range("E24").Value = "  +0.59%  "

$ws.Range("D25").Value = "2.172.39"
$ws.Range("E25").Value = "  +0.78%  "

$ws.Range("D26").Value = "6.659"
$ws.Range("E26").Value = "  +9.54%  "

$ws.Range("D27").Value = "155.63"
$ws.Range("E27").Value = "  +0.28%  "

$ws.Range("E28").Value = "  +1.12%  "

$ws.Range("E29").Value = "  +1.04%  "

$ws.Range("D30").Value = "120.53"
$ws.Range("E30").Value = "  +1.82%  "

$ws.Range("D31").Value = "1.015"
$ws.Range("E31").Value = "  -1.51%  "

$ws.Range("D32").Value = "0.09602"
$ws.Range("E32").Value = "  +0.96%  "

$ws.Range("D33").Value = "5.648"
$ws.Range("E33").Value = "  +4.65%  "

$ws.Range("D34").Value = "3.552"
$ws.Range("E34").Value = "  +0.43%  "

$ws.Range("E35").Value = "  -1.32%  "

$ws.Range("E36").Value = "  +1.18%  "

$ws.Range("D37").Value = "0.06101"
$ws.Range("E37").Value = "  +0.57%  "

$ws.Range("D38").Value = "1.179"
$ws.Range("E38").Value = "  +0.32%  "

$ws.Range("D39").Value = "8.074"
$ws.Range("E39").Value = "  +2.20%  "

$ws.Range("E42").Value = "  -0.07%  "

$ws.Range("E45").Value = "  -1.90%  "

$ws.Range("D46").Value = "12.40"
$ws.Range("E46").Value = "  +0.83%  "

$ws.Range("D47").Value = "0.5593"
$ws.Range("E47").Value = "  +1.27%  "

$ws.Range("D48").Value = "1.951"
$ws.Range("E48").Value = "  +1.53%  "

$ws.Range("D49").Value = "117.93"
$ws.Range("E49").Value = "  +3.96%  "

$ws.Range("D51").Value = "72.17"
$ws.Range("E51").Value = "  +0.05%  "

# --- Row 40/41 swap: TheSandbox <-> Aptos (Aptos now ranks above TheSandbox) ---
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "10.86"
$ws.Range("E40").Value = "  +5.84%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.5976"
$ws.Range("E41").Value = "  +1.67%  "

# --- Row 43/44 swap: WEMIXToken <-> RenderToken (RenderToken now ranks above WEMIXToken) ---
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "2.390"
$ws.Range("E43").Value = "  -1.37%  "

$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "1.247"
$ws.Range("E44").Value = "  -2.91%  "

